$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.226.79'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.856.28'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'0.6995"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = "'237.50"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.08024"
$ws.Range("E8").Value = '  +7.97%  '
$ws.Range("D9").Value = "'0.3041"
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").Value = "'23.21"
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = "'0.08178"
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '1.850.37'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = "'0.7142"
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = "'5.171"
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").Value = "'89.01"
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '29.224.38'
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = "'13.37"
$ws.Range("E17").Value = '  +2.82%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = "'5.751"
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = "'0.000007815"
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("D20").Value = "'235.58"
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '2.104.76'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D24").Value = "'7.409"
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").Value = "'161.70"
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").Value = "'8.956"
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").Value = "'0.1449"
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").Value = "'18.01"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = "'1.997"
$ws.Range("E29").Value = '  +3.26%  '
$ws.Range("D30").Value = "'1.433"
$ws.Range("E30").Value = '  +4.62%  '
$ws.Range("D31").Value = "'1.482"
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").Value = "'4.392"
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Value = "'4.047"
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").Value = "'0.05203"
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").Value = "'1.170"
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("D36").Value = "'0.7069"
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").Value = "'0.9986"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = "'2.663"
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("D39").Value = "'0.01845"
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = "'2.716"
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("D41").Value = '1.143.67'
$ws.Range("E41").Value = '  +8.47%  '
$ws.Range("D42").Value = "'0.9239"
$ws.Range("E42").Value = '  +2.57%  '
$ws.Range("D43").Value = "'5.938"
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").Value = "'0.4265"
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").Value = "'70.60"
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = "'102.87"
$ws.Range("E47").Value = '  +1.37%  '
$ws.Range("D48").Value = "'1.780"
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").Value = '2.000.88'
$ws.Range("E49").Value = '  +1.18%  '
$ws.Range("D50").Value = "'9.149"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = "'6.932"
$ws.Range("E51").Value = '  -1.71%  '
